$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-dimension:situacion-preferente -> iaest-measure:situacion-preferente (B2, shared with row1 "situacion-preferente")
$ws.Range("B2").Value = "iaest-measure:situacion-preferente"

# Row 2: municipio-nombre column (D2) and aragon column (F2) now reference refArea too
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "sdmx-dimension:refArea"

# Row 3: situacion-preferente column becomes "medida" instead of "dim"
$ws.Range("B3").Value = "medida"
# Row 3: municipio-nombre column becomes "dim" instead of "medida"
$ws.Range("D3").Value = "dim"

# Row 4: situacion-preferente column becomes "xsd:int" instead of "skos:Concept"
$ws.Range("B4").Value = "xsd:int"
# Row 4: municipio-nombre column gets new URI-Municipio value (was xsd:int)
$ws.Range("D4").Value = "URI-Municipio"
# Row 4: aragon column becomes "URI-Comunidad" instead of "skos:Concept"
$ws.Range("F4").Value = "URI-Comunidad"

# Row 5 (mapping-situacion-preferente.xlsx / mapping-aragon.xlsx) is removed entirely
$ws.Rows("5:5").Delete()
